$d = $word.ActiveDocument

# --- Step 1: merge runs split by proofErr markers (spell-check artifacts) ---
# These Find/Execute replacements span the proofErr boundaries and collapse
# the multiple runs into a single run, which also drops the now-orphaned
# <w:proofErr/> markers.

$rng1 = $d.Content
[void]$rng1.Find.Execute("disparo del cañon 1", $true, $false, $false, $false, $false, $true, 1, $false, "disparo del cañon 1", 2)

$rng2 = $d.Content
[void]$rng2.Find.Execute("(XD,YD) <= 0,025d(XO,YO)", $true, $false, $false, $false, $false, $true, 1, $false, "(XD,YD) <= 0,025d(XO,YO)", 2)

$rng3 = $d.Content
[void]$rng3.Find.Execute("posición Xo + Vx1*(T1+2).", $true, $false, $false, $false, $false, $true, 1, $false, "posición Xo + Vx1*(T1+2).", 2)

# --- Step 2: append the new paragraphs (the planning notes) before the
# trailing empty paragraph, using InsertXML so paragraph/run/tab structure
# comes out clean (no synthesized empty runs, explicit <w:tab/> elements). ---

$docEnd = $d.Content.End
$insertionPoint = $d.Range($docEnd, $docEnd)

$xmlFragment = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p/>
<w:p/>
<w:p/>
<w:p/>
<w:p/>
<w:p/>
<w:p/>
<w:p/>
<w:p/>
<w:p><w:r><w:t>Función</w:t></w:r><w:r><w:t xml:space="preserve"> calcular el movimiento con T+2</w:t></w:r></w:p>
<w:p><w:r><w:tab/><w:t>Recibe los parámetros</w:t></w:r></w:p>
<w:p><w:r><w:tab/><w:t>Retorna la posición final donde podría caer el disparo</w:t></w:r></w:p>
<w:p><w:r><w:t xml:space="preserve">Dado que los parámetros pueden varias </w:t></w:r><w:r><w:t>las funciones creadas deben de ser triviales y no excluir los distintos casos.</w:t></w:r></w:p>
<w:p><w:r><w:t>El programa debe recibir los parámetros de punto inicial del cañón ofensivo</w:t></w:r></w:p>
<w:p><w:r><w:t>Por ejemplo, posición inicial de el cañón, velocidad con la cual el disparo sale</w:t></w:r></w:p>
<w:p><w:r><w:t>Se debe calcular según la velocidad, la posición de el disparo con respecto a el tiempo</w:t></w:r></w:p>
<w:p><w:r><w:t>Se debe tener condicionales ya que el disparo ofensivo tiene mayor radio de explosión</w:t></w:r></w:p>
<w:p><w:r><w:t>Se debe hacer un condicional para que el disparo defensivo explote justo a una distancia de 0,025 de el disparo ofensivo para que este no lo detecte.</w:t></w:r></w:p>
<w:p><w:r><w:t>Se deben hacer los distintos condicionales para las diferentes situaciones planteadas</w:t></w:r></w:p>
<w:p><w:r><w:t>También hay que añadirle a el disparo defensivo que es el de salida un retraso de 0,5 segundos</w:t></w:r></w:p>
<w:p><w:r><w:t>Dado esto se podría hacer una función que dada “</w:t></w:r><w:r><w:t>Función</w:t></w:r><w:r><w:t xml:space="preserve"> calcular movimiento con T+2</w:t></w:r><w:r><w:t>” sirva para la exactitud de los parámetros de salida</w:t></w:r></w:p>
<w:p><w:r><w:t>Se puede hacer una secuencia dado que puede generar al menos 3 disparos en cada una de las situaciones.</w:t></w:r></w:p>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

[void]$insertionPoint.InsertXML($xmlFragment)

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
